# Fruta / hortaliza, semanal
# Update the Chirimoya (La Palmera de La Serena) weekly price table:
#  - remove the old "Extra (doble especial)" row for Provincia de Limari / 10-09-2021
#  - insert a new weekly block (Especial / Primera / Segunda) for
#    Provincia de Limari dated 09-11-2021 (serial 44509)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the obsolete row 84 (Extra (doble especial), Limari, 44449) - everything
#    below shifts up by one row.
$ws.Rows.Item(84).Delete()

# 2) After the deletion, the old row 96 (Especial / Provincia del Elqui / 44491)
#    is now row 95. Insert three fresh rows above it to hold the new weekly data.
$ws.Rows.Item(95).Insert()
$ws.Rows.Item(95).Insert()
$ws.Rows.Item(95).Insert()

# 3) Populate the three new rows (95-97) with the new Limari week (44509).
$ws.Cells.Item(95, 1).Value = 8
$ws.Cells.Item(95, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 44509
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100107
$ws.Cells.Item(95, 8).Value = "Otros"
$ws.Cells.Item(95, 9).Value = 100107002
$ws.Cells.Item(95, 10).Value = "Chirimoya"
$ws.Cells.Item(95, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(95, 12).Value = "Especial"
$ws.Cells.Item(95, 13).Value = 500
$ws.Cells.Item(95, 14).Value = 1800
$ws.Cells.Item(95, 15).Value = 1900
$ws.Cells.Item(95, 16).Value = 1850
$ws.Cells.Item(95, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(95, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(95, 19).Value = 1850
$ws.Cells.Item(95, 20).Value = 1

$ws.Cells.Item(96, 1).Value = 8
$ws.Cells.Item(96, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(96, 3).Value = "Coquimbo"
$ws.Cells.Item(96, 4).Value = 44509
$ws.Cells.Item(96, 5).Value = 4
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100107
$ws.Cells.Item(96, 8).Value = "Otros"
$ws.Cells.Item(96, 9).Value = 100107002
$ws.Cells.Item(96, 10).Value = "Chirimoya"
$ws.Cells.Item(96, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 400
$ws.Cells.Item(96, 14).Value = 1500
$ws.Cells.Item(96, 15).Value = 1600
$ws.Cells.Item(96, 16).Value = 1550
$ws.Cells.Item(96, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(96, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(96, 19).Value = 1550
$ws.Cells.Item(96, 20).Value = 1

$ws.Cells.Item(97, 1).Value = 8
$ws.Cells.Item(97, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(97, 3).Value = "Coquimbo"
$ws.Cells.Item(97, 4).Value = 44509
$ws.Cells.Item(97, 5).Value = 4
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100107
$ws.Cells.Item(97, 8).Value = "Otros"
$ws.Cells.Item(97, 9).Value = 100107002
$ws.Cells.Item(97, 10).Value = "Chirimoya"
$ws.Cells.Item(97, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(97, 12).Value = "Segunda"
$ws.Cells.Item(97, 13).Value = 280
$ws.Cells.Item(97, 14).Value = 1300
$ws.Cells.Item(97, 15).Value = 1400
$ws.Cells.Item(97, 16).Value = 1350
$ws.Cells.Item(97, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(97, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(97, 19).Value = 1350
$ws.Cells.Item(97, 20).Value = 1
